# Auto-generated Excel COM-interop edit script
# Commit: Add data for 2025-04-23
# Applies 150 cell value updates across 40 worksheets (daily crime-count increments
# for 2025 year-to-date figures, plus a handful of minor prior-year corrections).

$wb = $excel.ActiveWorkbook

# --- Sheet: Citywide Totals ---
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Cells.Item(2, 12).Value = 1821   # L2: 1804 -> 1821
$ws.Cells.Item(3, 12).Value = 1851   # L3: 1831 -> 1851
$ws.Cells.Item(4, 3).Value = 1867   # C4: 1866 -> 1867
$ws.Cells.Item(4, 9).Value = 1834   # I4: 1835 -> 1834
$ws.Cells.Item(4, 11).Value = 1756   # K4: 1755 -> 1756
$ws.Cells.Item(4, 12).Value = 523   # L4: 515 -> 523
$ws.Cells.Item(6, 12).Value = 1713   # L6: 1695 -> 1713
$ws.Cells.Item(7, 3).Value = 28411   # C7: 28410 -> 28411
$ws.Cells.Item(7, 9).Value = 26301   # I7: 26302 -> 26301
$ws.Cells.Item(7, 11).Value = 27547   # K7: 27546 -> 27547
$ws.Cells.Item(7, 12).Value = 6019   # L7: 5956 -> 6019

# --- Sheet: Logan Square ---
$ws = $wb.Worksheets.Item("Logan Square")
$ws.Cells.Item(6, 12).Value = 26   # L6: 25 -> 26
$ws.Cells.Item(7, 12).Value = 73   # L7: 72 -> 73

# --- Sheet: Austin ---
$ws = $wb.Worksheets.Item("Austin")
$ws.Cells.Item(3, 12).Value = 123   # L3: 122 -> 123
$ws.Cells.Item(6, 12).Value = 104   # L6: 103 -> 104
$ws.Cells.Item(7, 12).Value = 370   # L7: 368 -> 370

# --- Sheet: Garfield Park ---
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Cells.Item(2, 12).Value = 62   # L2: 60 -> 62
$ws.Cells.Item(3, 12).Value = 91   # L3: 90 -> 91
$ws.Cells.Item(7, 12).Value = 264   # L7: 261 -> 264

# --- Sheet: Grand Crossing ---
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Cells.Item(2, 12).Value = 69   # L2: 68 -> 69
$ws.Cells.Item(3, 12).Value = 59   # L3: 58 -> 59
$ws.Cells.Item(7, 12).Value = 216   # L7: 214 -> 216

# --- Sheet: New City ---
$ws = $wb.Worksheets.Item("New City")
$ws.Cells.Item(3, 12).Value = 35   # L3: 34 -> 35
$ws.Cells.Item(7, 12).Value = 116   # L7: 115 -> 116

# --- Sheet: Woodlawn ---
$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Cells.Item(2, 12).Value = 25   # L2: 24 -> 25
$ws.Cells.Item(7, 12).Value = 95   # L7: 94 -> 95

# --- Sheet: By Neighborhood ---
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Cells.Item(7, 12).Value = 202   # L7: 199 -> 202
$ws.Cells.Item(8, 12).Value = 370   # L8: 368 -> 370
$ws.Cells.Item(9, 12).Value = 38   # L9: 36 -> 38
$ws.Cells.Item(10, 12).Value = 40   # L10: 38 -> 40
$ws.Cells.Item(15, 12).Value = 40   # L15: 38 -> 40
$ws.Cells.Item(18, 12).Value = 46   # L18: 44 -> 46
$ws.Cells.Item(20, 12).Value = 161   # L20: 160 -> 161
$ws.Cells.Item(21, 12).Value = 19   # L21: 18 -> 19
$ws.Cells.Item(23, 12).Value = 63   # L23: 62 -> 63
$ws.Cells.Item(25, 12).Value = 32   # L25: 31 -> 32
$ws.Cells.Item(27, 11).Value = 262   # K27: 261 -> 262
$ws.Cells.Item(27, 12).Value = 64   # L27: 63 -> 64
$ws.Cells.Item(29, 12).Value = 307   # L29: 306 -> 307
$ws.Cells.Item(33, 12).Value = 264   # L33: 261 -> 264
$ws.Cells.Item(37, 12).Value = 216   # L37: 214 -> 216
$ws.Cells.Item(42, 12).Value = 188   # L42: 186 -> 188
$ws.Cells.Item(43, 12).Value = 49   # L43: 47 -> 49
$ws.Cells.Item(48, 12).Value = 90   # L48: 89 -> 90
$ws.Cells.Item(50, 12).Value = 36   # L50: 35 -> 36
$ws.Cells.Item(51, 12).Value = 73   # L51: 72 -> 73
$ws.Cells.Item(52, 12).Value = 128   # L52: 126 -> 128
$ws.Cells.Item(53, 12).Value = 73   # L53: 72 -> 73
$ws.Cells.Item(54, 12).Value = 130   # L54: 128 -> 130
$ws.Cells.Item(59, 12).Value = 10   # L59: 9 -> 10
$ws.Cells.Item(60, 12).Value = 34   # L60: 33 -> 34
$ws.Cells.Item(63, 3).Value = 292   # C63: 291 -> 292
$ws.Cells.Item(63, 9).Value = 256   # I63: 257 -> 256
$ws.Cells.Item(63, 11).Value = 87   # K63: 89 -> 87
$ws.Cells.Item(63, 12).Value = 22   # L63: 19 -> 22
$ws.Cells.Item(65, 12).Value = 116   # L65: 115 -> 116
$ws.Cells.Item(67, 12).Value = 210   # L67: 207 -> 210
$ws.Cells.Item(70, 11).Value = 50   # K70: 49 -> 50
$ws.Cells.Item(72, 12).Value = 25   # L72: 24 -> 25
$ws.Cells.Item(79, 11).Value = 671   # K79: 670 -> 671
$ws.Cells.Item(84, 12).Value = 61   # L84: 60 -> 61
$ws.Cells.Item(85, 12).Value = 316   # L85: 308 -> 316
$ws.Cells.Item(86, 12).Value = 43   # L86: 41 -> 43
$ws.Cells.Item(89, 12).Value = 73   # L89: 72 -> 73
$ws.Cells.Item(90, 12).Value = 59   # L90: 58 -> 59
$ws.Cells.Item(91, 12).Value = 78   # L91: 77 -> 78
$ws.Cells.Item(93, 12).Value = 31   # L93: 30 -> 31
$ws.Cells.Item(98, 12).Value = 45   # L98: 44 -> 45
$ws.Cells.Item(99, 12).Value = 95   # L99: 94 -> 95
$ws.Cells.Item(101, 3).Value = 28411   # C101: 28410 -> 28411
$ws.Cells.Item(101, 9).Value = 26301   # I101: 26302 -> 26301
$ws.Cells.Item(101, 11).Value = 27547   # K101: 27546 -> 27547
$ws.Cells.Item(101, 12).Value = 6019   # L101: 5956 -> 6019

# --- Sheet: North Lawndale ---
$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Cells.Item(2, 12).Value = 61   # L2: 60 -> 61
$ws.Cells.Item(3, 12).Value = 65   # L3: 64 -> 65
$ws.Cells.Item(6, 12).Value = 58   # L6: 57 -> 58
$ws.Cells.Item(7, 12).Value = 210   # L7: 207 -> 210

# --- Sheet: South Deering ---
$ws = $wb.Worksheets.Item("South Deering")
$ws.Cells.Item(2, 12).Value = 24   # L2: 23 -> 24
$ws.Cells.Item(7, 12).Value = 61   # L7: 60 -> 61

# --- Sheet: Loop ---
$ws = $wb.Worksheets.Item("Loop")
$ws.Cells.Item(2, 12).Value = 33   # L2: 32 -> 33
$ws.Cells.Item(6, 12).Value = 64   # L6: 63 -> 64
$ws.Cells.Item(7, 12).Value = 130   # L7: 128 -> 130

# --- Sheet: Englewood ---
$ws = $wb.Worksheets.Item("Englewood")
$ws.Cells.Item(2, 12).Value = 98   # L2: 97 -> 98
$ws.Cells.Item(7, 12).Value = 307   # L7: 306 -> 307

# --- Sheet: Lake View ---
$ws = $wb.Worksheets.Item("Lake View")
$ws.Cells.Item(3, 12).Value = 19   # L3: 18 -> 19
$ws.Cells.Item(7, 12).Value = 90   # L7: 89 -> 90

# --- Sheet: Humboldt Park ---
$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Cells.Item(3, 12).Value = 51   # L3: 49 -> 51
$ws.Cells.Item(7, 12).Value = 188   # L7: 186 -> 188

# --- Sheet: Avondale ---
$ws = $wb.Worksheets.Item("Avondale")
$ws.Cells.Item(2, 12).Value = 18   # L2: 17 -> 18
$ws.Cells.Item(6, 12).Value = 12   # L6: 11 -> 12
$ws.Cells.Item(7, 12).Value = 40   # L7: 38 -> 40

# --- Sheet: Douglas ---
$ws = $wb.Worksheets.Item("Douglas")
$ws.Cells.Item(6, 12).Value = 11   # L6: 10 -> 11
$ws.Cells.Item(7, 12).Value = 63   # L7: 62 -> 63

# --- Sheet: Washington Park ---
$ws = $wb.Worksheets.Item("Washington Park")
$ws.Cells.Item(3, 12).Value = 26   # L3: 25 -> 26
$ws.Cells.Item(7, 12).Value = 78   # L7: 77 -> 78
$ws.Cells.Item(6, 12).Value = 13   # L6: 12 -> 13

# --- Sheet: Chinatown ---
$ws = $wb.Worksheets.Item("Chinatown")
$ws.Cells.Item(7, 12).Value = 19   # L7: 18 -> 19

# --- Sheet: Roseland ---
$ws = $wb.Worksheets.Item("Roseland")
$ws.Cells.Item(4, 11).Value = 43   # K4: 42 -> 43
$ws.Cells.Item(7, 11).Value = 671   # K7: 670 -> 671

# --- Sheet: Chicago Lawn ---
$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Cells.Item(3, 12).Value = 46   # L3: 45 -> 46
$ws.Cells.Item(7, 12).Value = 161   # L7: 160 -> 161

# --- Sheet: Calumet Heights ---
$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Cells.Item(2, 12).Value = 17   # L2: 16 -> 17
$ws.Cells.Item(6, 12).Value = 6   # L6: 5 -> 6
$ws.Cells.Item(7, 12).Value = 46   # L7: 44 -> 46

# --- Sheet: West Lawn ---
$ws = $wb.Worksheets.Item("West Lawn")
$ws.Cells.Item(3, 12).Value = 7   # L3: 6 -> 7
$ws.Cells.Item(7, 12).Value = 31   # L7: 30 -> 31

# --- Sheet: Auburn Gresham ---
$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Cells.Item(3, 12).Value = 64   # L3: 61 -> 64
$ws.Cells.Item(7, 12).Value = 202   # L7: 199 -> 202

# --- Sheet: East Side ---
$ws = $wb.Worksheets.Item("East Side")
$ws.Cells.Item(6, 12).Value = 6   # L6: 5 -> 6
$ws.Cells.Item(7, 12).Value = 32   # L7: 31 -> 32

# --- Sheet: Brighton Park ---
$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Cells.Item(3, 12).Value = 17   # L3: 15 -> 17
$ws.Cells.Item(7, 12).Value = 40   # L7: 38 -> 40

# --- Sheet: Wicker Park ---
$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Cells.Item(6, 12).Value = 25   # L6: 24 -> 25
$ws.Cells.Item(7, 12).Value = 45   # L7: 44 -> 45

# --- Sheet: Lincoln Square ---
$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Cells.Item(6, 12).Value = 8   # L6: 7 -> 8
$ws.Cells.Item(7, 12).Value = 36   # L7: 35 -> 36

# --- Sheet: Avalon Park ---
$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Cells.Item(2, 12).Value = 9   # L2: 8 -> 9
$ws.Cells.Item(4, 12).Value = 2   # L4: 1 -> 2
$ws.Cells.Item(7, 12).Value = 38   # L7: 36 -> 38

# --- Sheet: Montclare ---
$ws = $wb.Worksheets.Item("Montclare")
$ws.Cells.Item(3, 12).Value = 5   # L3: 4 -> 5
$ws.Cells.Item(7, 12).Value = 10   # L7: 9 -> 10

# --- Sheet: O'Hare ---
$ws = $wb.Worksheets.Item("O'Hare")
$ws.Cells.Item(6, 11).Value = 7   # K6: 6 -> 7
$ws.Cells.Item(7, 11).Value = 50   # K7: 49 -> 50

# --- Sheet: Uptown ---
$ws = $wb.Worksheets.Item("Uptown")
$ws.Cells.Item(6, 12).Value = 14   # L6: 13 -> 14
$ws.Cells.Item(7, 12).Value = 73   # L7: 72 -> 73

# --- Sheet: Edgewater ---
$ws = $wb.Worksheets.Item("Edgewater")
$ws.Cells.Item(3, 12).Value = 25   # L3: 24 -> 25
$ws.Cells.Item(4, 11).Value = 34   # K4: 33 -> 34
$ws.Cells.Item(7, 11).Value = 262   # K7: 261 -> 262
$ws.Cells.Item(7, 12).Value = 64   # L7: 63 -> 64

# --- Sheet: Streeterville ---
$ws = $wb.Worksheets.Item("Streeterville")
$ws.Cells.Item(3, 12).Value = 7   # L3: 6 -> 7
$ws.Cells.Item(4, 12).Value = 24   # L4: 23 -> 24
$ws.Cells.Item(7, 12).Value = 43   # L7: 41 -> 43

# --- Sheet: Washington Heights ---
$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Cells.Item(2, 12).Value = 24   # L2: 23 -> 24
$ws.Cells.Item(7, 12).Value = 59   # L7: 58 -> 59

# --- Sheet: Little Italy, UIC ---
$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Cells.Item(6, 12).Value = 21   # L6: 20 -> 21
$ws.Cells.Item(7, 12).Value = 73   # L7: 72 -> 73

# --- Sheet: Morgan Park ---
$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Cells.Item(6, 12).Value = 10   # L6: 9 -> 10
$ws.Cells.Item(7, 12).Value = 34   # L7: 33 -> 34

# --- Sheet: Hyde Park ---
$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Cells.Item(4, 12).Value = 11   # L4: 10 -> 11
$ws.Cells.Item(6, 12).Value = 16   # L6: 15 -> 16
$ws.Cells.Item(7, 12).Value = 49   # L7: 47 -> 49

# --- Sheet: South Shore ---
$ws = $wb.Worksheets.Item("South Shore")
$ws.Cells.Item(2, 12).Value = 95   # L2: 92 -> 95
$ws.Cells.Item(3, 12).Value = 130   # L3: 128 -> 130
$ws.Cells.Item(6, 12).Value = 57   # L6: 54 -> 57
$ws.Cells.Item(7, 12).Value = 316   # L7: 308 -> 316

# --- Sheet: Old Town ---
$ws = $wb.Worksheets.Item("Old Town")
$ws.Cells.Item(4, 12).Value = 4   # L4: 3 -> 4
$ws.Cells.Item(7, 12).Value = 25   # L7: 24 -> 25

# --- Sheet: Little Village ---
$ws = $wb.Worksheets.Item("Little Village")
$ws.Cells.Item(2, 12).Value = 42   # L2: 41 -> 42
$ws.Cells.Item(3, 12).Value = 37   # L3: 36 -> 37
$ws.Cells.Item(7, 12).Value = 128   # L7: 126 -> 128
